$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 756
$wsExhibit.Range("F3").Value = 13
$wsExhibit.Range("F4").Value = 50
$wsExhibit.Range("F5").Value = 23
$wsExhibit.Range("F7").Value = 3424
$wsExhibit.Range("F9").Value = 4098
$wsExhibit.Range("F11").Value = 1024
$wsExhibit.Range("F12").Value = 44

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 756
$wsAll.Range("F3").Value = 13
$wsAll.Range("F4").Value = 50
$wsAll.Range("F5").Value = 23
$wsAll.Range("F8").Value = 3424
$wsAll.Range("F10").Value = 4098
$wsAll.Range("F12").Value = 1024
$wsAll.Range("F13").Value = 44
